# Regenerate save_data to use K instead of Strike#, recalc std/mean, and write
# the updated s_vals (K column, column G) back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (s_vals) values keyed by row number (row 1 is the header row).
$newK = @{
    2  = 4
    3  = 6
    4  = 1
    5  = 6
    6  = 5
    7  = 12
    8  = 7
    9  = 8
    10 = 7
    11 = 0
    12 = 2
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 0
    18 = 1
    19 = 4
    20 = 0
    21 = 2
    22 = 3
    23 = 2
    24 = 0
    25 = 2
    26 = 3
    27 = 5
    28 = 7
    29 = 3
    30 = 2
    31 = 5
    32 = 3
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
